$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "29.711.66"
$ws.Range("E2").Value = "  -3.43%  "
$ws.Range("D3").Value = "2.096.08"
$ws.Range("E3").Value = "  -2.56%  "
Set-TextValue $ws.Range("D4") "1.009"
$ws.Range("E4").Value = "  -0.30%  "
Set-TextValue $ws.Range("D5") "344.78"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("E6").Value = "  -0.23%  "
Set-TextValue $ws.Range("D7") "0.5146"
$ws.Range("E7").Value = "  -2.71%  "
Set-TextValue $ws.Range("D8") "0.4406"
$ws.Range("E8").Value = "  -3.78%  "
Set-TextValue $ws.Range("D9") "52.69"
$ws.Range("E9").Value = "  -3.03%  "
Set-TextValue $ws.Range("D10") "0.09274"
$ws.Range("E10").Value = "  +0.68%  "
Set-TextValue $ws.Range("D11") "1.171"
$ws.Range("E11").Value = "  -1.43%  "
Set-TextValue $ws.Range("D12") "24.85"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "2.105.20"
$ws.Range("E13").Value = "  -2.03%  "
Set-TextValue $ws.Range("D14") "8.294"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("E15").Value = "  -2.51%  "
$ws.Range("E16").Value = "  -2.81%  "
Set-TextValue $ws.Range("D17") "0.00001153"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("E18").Value = "  -0.26%  "
Set-TextValue $ws.Range("D19") "20.87"
$ws.Range("E19").Value = "  +6.11%  "
Set-TextValue $ws.Range("D20") "0.06629"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  -0.23%  "
Set-TextValue $ws.Range("D22") "6.198"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "29.768.52"
$ws.Range("E23").Value = "  -3.54%  "
Set-TextValue $ws.Range("D24") "12.60"
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("D26").Value = "2.352.76"
$ws.Range("E26").Value = "  -1.32%  "
Set-TextValue $ws.Range("D27") "21.89"
$ws.Range("E27").Value = "  -3.32%  "
Set-TextValue $ws.Range("D28") "2.525"
$ws.Range("E28").Value = "  -4.26%  "
Set-TextValue $ws.Range("D29") "161.95"
$ws.Range("E29").Value = "  -2.01%  "
Set-TextValue $ws.Range("D30") "132.80"
$ws.Range("E30").Value = "  -3.12%  "
Set-TextValue $ws.Range("D31") "1.132"
$ws.Range("E31").Value = "  -7.85%  "
Set-TextValue $ws.Range("D32") "0.1051"
$ws.Range("E32").Value = "  -3.10%  "
Set-TextValue $ws.Range("D33") "1.656"
$ws.Range("E33").Value = "  -1.37%  "
Set-TextValue $ws.Range("D34") "6.178"
$ws.Range("E34").Value = "  -3.76%  "
Set-TextValue $ws.Range("D35") "3.943"
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D36") "6.095"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D37") "10.48"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -3.58%  "
Set-TextValue $ws.Range("D39") "0.06726"
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("E40").Value = "  -2.09%  "
Set-TextValue $ws.Range("D41") "0.6863"
$ws.Range("E41").Value = "  -2.07%  "
Set-TextValue $ws.Range("D42") "0.2226"
$ws.Range("E42").Value = "  -5.03%  "
Set-TextValue $ws.Range("D43") "1.300"
$ws.Range("E43").Value = "  +1.76%  "
Set-TextValue $ws.Range("D44") "0.6639"
Set-TextValue $ws.Range("D45") "14.25"
$ws.Range("E45").Value = "  -3.81%  "
Set-TextValue $ws.Range("D46") "2.317"
$ws.Range("E47").Value = "  -3.80%  "
Set-TextValue $ws.Range("D48") "0.00000000348"
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D50") "82.39"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue $ws.Range("D51") "0.3335"
$ws.Range("E51").Value = "  -0.62%  "
